$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Fix the mis-nested runs in the "Front End Web Developer" contact-info
#    paragraph. The source document has an empty wrapper <w:r> (sz 22 Arial)
#    that illegally contains four more <w:r> elements. Re-insert the
#    paragraph's XML with those runs flattened as siblings (no wrapper run)
#    so Google Docs / LibreOffice parse it correctly.
# ---------------------------------------------------------------------------
$contactPara = $d.Paragraphs(2)

$paraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:color w:val="444444"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
      <w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/>
    </w:rPr>
    <w:t xml:space="preserve">Front End Web Developer</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:color w:val="D4D4D4"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/>
    </w:rPr>
    <w:t xml:space="preserve">    |    </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:i/>
      <w:iCs/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
      <w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/>
    </w:rPr>
    <w:t xml:space="preserve">+1 (647) 401-1468</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
      <w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/>
    </w:rPr>
    <w:t xml:space="preserve"> &#8226; </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:i/>
      <w:iCs/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
      <w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/>
    </w:rPr>
    <w:t xml:space="preserve">oleg@olegrybin.com</w:t>
  </w:r>
</w:p>
'@

$contactPara.Range.InsertXML($paraXml)

# ---------------------------------------------------------------------------
# 2) Make the "Key skills" table compatible with Google Docs / Open Office:
#    switch the table width from a percentage (pct 100%) to auto (w="100"),
#    and switch the three grid columns / cell widths from percentage-based
#    (pct 33%, grid placeholder 100) to fixed dxa widths of 3500 twips each.
# ---------------------------------------------------------------------------
$tbl = $d.Tables(1)

# wdPreferredWidthAuto = 1 -> renders as w:type="auto"; width in points,
# 5pt * 20 = 100 twips, matching <w:tblW w:type="auto" w:w="100"/>.
$tbl.PreferredWidthType = 1
$tbl.PreferredWidth = 5

# 175pt * 20 = 3500 twips for each of the three columns.
for ($i = 1; $i -le $tbl.Columns.Count; $i++) {
    $tbl.Columns($i).Width = 175
}

# wdPreferredWidthPoints = 3 -> renders as w:type="dxa".
for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
        $cell = $tbl.Cell($r, $c)
        $cell.PreferredWidthType = 3
        $cell.PreferredWidth = 175
    }
}

Write-Host "Edit complete"
